$d = $word.ActiveDocument

# 1) " [1/5]" -> "   [1/5]"  (after the hyperlink, in the first paragraph)
# Search just for "[1/5]" (not including the leading space) so the match
# doesn't start exactly on the hyperlink-run boundary, which would make the
# replacement inherit the hyperlink's character style.
$r1 = $d.Content
$r1.Find.Execute("[1/5]", $false, $false, $false, $false, $false, $true, 1, $false, "  [1/5]", 2)

# 2) "...especially on very niche pages. " -> "...especially when analyzing very niche topics. "
$r2 = $d.Content
$r2.Find.Execute("especially on very niche pages.", $false, $false, $false, $false, $false, $true, 1, $false, "especially when analyzing very niche topics.", 2)

# 3) " conitnued interest in this, I'll keep..." -> " continued interest in this little project, I'll keep..."
$r3 = $d.Content
$r3.Find.Execute(" conitnued interest in this,", $false, $false, $false, $false, $false, $true, 1, $false, " continued interest in this little project,", 2)
